# Add 2022-Q4 data:
#  - Insert a new "2022-Q4" worksheet right after "总计" (pushing the
#    existing quarter sheets down / renaming them to the next quarter slot).
#  - Populate the new sheet with the 2022-Q4 fund holdings.
#  - Insert a new summary row on "总计" for 2022-Q4 and keep the older
#    rows (now shifted down one row) consistent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet immediately after "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# The sheet collection got rebuilt by Add(), so re-fetch the sheet we
# want to use as a formatting template (previously "2022-Q3").
$q3 = $wb.Worksheets.Item("2022-Q3")

# Copy over formatting (header styling, index-column styling, column
# widths, etc.) from the 2022-Q3 sheet so the new sheet matches the
# existing look & feel.
$q3.Range("A1:H5").Copy()
$q4.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Fill in the 2022-Q4 fund holdings data.
# ---------------------------------------------------------------------
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "'014317"
$q4.Cells.Item(2,3).Value = "广发价值领航一年持有混合A"
$q4.Cells.Item(2,4).Value = "'2.40"
$q4.Cells.Item(2,5).Value = "'93.40"
$q4.Cells.Item(2,6).Value = "'4.25"
$q4.Cells.Item(2,7).Value = "'0.1020"
$q4.Cells.Item(2,8).Value = 10

$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "'011351"
$q4.Cells.Item(3,3).Value = "金鹰年年邮益一年持有期混合A"
$q4.Cells.Item(3,4).Value = "'3.04"
$q4.Cells.Item(3,5).Value = "'39.17"
$q4.Cells.Item(3,6).Value = "'0.93"
$q4.Cells.Item(3,7).Value = "'0.0283"
$q4.Cells.Item(3,8).Value = 5

$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "'014318"
$q4.Cells.Item(4,3).Value = "广发价值领航一年持有混合C"
$q4.Cells.Item(4,4).Value = "'0.65"
$q4.Cells.Item(4,5).Value = "'93.40"
$q4.Cells.Item(4,6).Value = "'4.25"
$q4.Cells.Item(4,7).Value = "'0.0276"
$q4.Cells.Item(4,8).Value = 10

$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "'011352"
$q4.Cells.Item(5,3).Value = "金鹰年年邮益一年持有期混合C"
$q4.Cells.Item(5,4).Value = "'0.23"
$q4.Cells.Item(5,5).Value = "'39.17"
$q4.Cells.Item(5,6).Value = "'0.93"
$q4.Cells.Item(5,7).Value = "'0.0021"
$q4.Cells.Item(5,8).Value = 5

# ---------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    above the existing 2022-Q3 row, shifting the rest down.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q4"
$totals.Cells.Item(2,3).Value = 4
$totals.Cells.Item(2,4).Value = 0.16

# Carry over formatting for the newly inserted row from the row below it
# (which holds what used to be row 2, i.e. the same original styling).
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$totals.Range("B3:D3").Copy()
$totals.Range("B2:D2").PasteSpecial(-4122)  # xlPasteFormats

# Renumber the index column (A) for all the data rows so it stays a
# simple 0-based sequence: 0,1,2,3,4.
for ($r = 3; $r -le 6; $r++) {
    $totals.Cells.Item($r, 1).Value = $r - 2
}
